$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before the current row 286, shifting the existing
# rows 286-304 down to 288-306 (same as the author's edit: two new weekly
# price records were added in the middle of the table).
$ws.Rows.Item(286).EntireRow.Insert()
$ws.Rows.Item(286).EntireRow.Insert()

# --- New row 286 ---
$ws.Cells.Item(286,1).Value  = 4
$ws.Cells.Item(286,2).Value  = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(286,3).Value  = "Los Lagos"
$ws.Cells.Item(286,4).Value  = (Get-Date -Year 2022 -Month 2 -Day 18 -Hour 0 -Minute 0 -Second 0)
$ws.Cells.Item(286,5).Value  = 10
$ws.Cells.Item(286,6).Value  = 100114013
$ws.Cells.Item(286,7).Value  = "Zanahoria"
$ws.Cells.Item(286,8).Value  = "Sin especificar"
$ws.Cells.Item(286,9).Value  = "Primera"
$ws.Cells.Item(286,10).Value = 350
$ws.Cells.Item(286,11).Value = 12000
$ws.Cells.Item(286,12).Value = 12000
$ws.Cells.Item(286,13).Value = 12000
$ws.Cells.Item(286,14).Value = "`$/saco 20 kilos"
$ws.Cells.Item(286,15).Value = "Chillán"
$ws.Cells.Item(286,16).Value = 600
$ws.Cells.Item(286,17).Value = 20
$ws.Cells.Item(286,18).Value = "Hortaliza"

# --- New row 287 ---
$ws.Cells.Item(287,1).Value  = 4
$ws.Cells.Item(287,2).Value  = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(287,3).Value  = "Los Lagos"
$ws.Cells.Item(287,4).Value  = (Get-Date -Year 2022 -Month 2 -Day 18 -Hour 0 -Minute 0 -Second 0)
$ws.Cells.Item(287,5).Value  = 10
$ws.Cells.Item(287,6).Value  = 100114013
$ws.Cells.Item(287,7).Value  = "Zanahoria"
$ws.Cells.Item(287,8).Value  = "Sin especificar"
$ws.Cells.Item(287,9).Value  = "Primera"
$ws.Cells.Item(287,10).Value = 500
$ws.Cells.Item(287,11).Value = 10000
$ws.Cells.Item(287,12).Value = 10000
$ws.Cells.Item(287,13).Value = 10000
$ws.Cells.Item(287,14).Value = "`$/saco 20 kilos"
$ws.Cells.Item(287,15).Value = "Provincia de Llanquihue"
$ws.Cells.Item(287,16).Value = 500
$ws.Cells.Item(287,17).Value = 20
$ws.Cells.Item(287,18).Value = "Hortaliza"
